$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.133.54"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "3.144.78"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.71"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.40"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +9.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +5.29%  "
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "3.679.85"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.03"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").Value = "58.179.68"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("E17").Value = "  +6.10%  "
$ws.Range("D18").Value = "3.139.20"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +4.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.12"
$ws.Range("E21").Value = "  +7.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.45"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.517"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  +12.45%  "
$ws.Range("D29").Value = "0.0₃0881"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  +7.32%  "
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("E33").Value = "  +7.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.54"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.36"
$ws.Range("E37").Value = "  +9.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.47"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +7.94%  "
$ws.Range("D40").Value = "2.645.81"
$ws.Range("E40").Value = "  +9.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0679"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").Value = "  +5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.54"
$ws.Range("E43").Value = "  +5.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.701"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +12.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.23"
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.980"
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.27"
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.753"
$ws.Range("E51").Value = "  +0.50%  "
